$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 15
$ws.Range("N15").Value = 1

# Row 19
$ws.Range("M19").Value = 2
$ws.Range("N19").Value = 2

# Row 31
$ws.Range("M31").Value = 3
$ws.Range("N31").Value = 3

# Row 32
$ws.Range("M32").Value = 1

# Row 34
$ws.Range("M34").Value = 1
$ws.Range("N34").Value = 1

# Row 38
$ws.Range("M38").Value = 2
$ws.Range("N38").Value = 2

# Row 42
$ws.Range("M42").Value = 4

# Row 45
$ws.Range("M45").Value = 3
$ws.Range("N45").Value = 3
